# Equip compose system optimisation:
# The "建造" (Build) icon entry (Id=10, Icon="MainIcon10") is no longer
# needed, so remove its entire row from the MainIcon table. Removing the
# worksheet row also shrinks the backing Excel Table / AutoFilter range
# and the sheet dimension automatically, and every row below shifts up
# by one (their values/styles travel with them).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 holds: A="10", B="建造", C="查看我的建造目录", N="MainIcon10"
$ws.Rows(11).Delete()
